$wb = $excel.ActiveWorkbook

# Update the "Status" value from "Ready for handoff" to "In Translation"
# on the zh-cn and de-de sheets (column C, row 2), and on the Overview
# sheet (columns E and F, row 2) which mirror the same status per-locale.

$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = "In Translation"

$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = "In Translation"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

# Re-fit the affected columns now that the text is shorter.
$newWidth = 13.4101845877511
$zh.Columns.Item(3).ColumnWidth = $newWidth
$de.Columns.Item(3).ColumnWidth = $newWidth
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth
